$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) '67.078.46'
Set-TextValue $ws.Cells.Item(2, 5) '  +1.17%  '

Set-TextValue $ws.Cells.Item(3, 4) '3.117.32'
Set-TextValue $ws.Cells.Item(3, 5) '  +2.79%  '

Set-TextValue $ws.Cells.Item(4, 5) '  +0.06%  '

Set-TextValue $ws.Cells.Item(5, 4) '580.02'
Set-TextValue $ws.Cells.Item(5, 5) '  +0.47%  '

Set-TextValue $ws.Cells.Item(6, 4) '174.23'
Set-TextValue $ws.Cells.Item(6, 5) '  +3.62%  '

Set-TextValue $ws.Cells.Item(7, 5) '  +0.00%  '

Set-TextValue $ws.Cells.Item(8, 4) '3.111.54'
Set-TextValue $ws.Cells.Item(8, 5) '  +2.71%  '

Set-TextValue $ws.Cells.Item(9, 5) '  +0.77%  '

Set-TextValue $ws.Cells.Item(10, 4) '6.46'
Set-TextValue $ws.Cells.Item(10, 5) '  -3.12%  '

Set-TextValue $ws.Cells.Item(11, 5) '  +1.24%  '

Set-TextValue $ws.Cells.Item(12, 5) '  +0.40%  '

Set-TextValue $ws.Cells.Item(13, 5) '  +0.33%  '

Set-TextValue $ws.Cells.Item(14, 4) '37.37'
Set-TextValue $ws.Cells.Item(14, 5) '  +2.73%  '

Set-TextValue $ws.Cells.Item(16, 4) '3.629.37'
Set-TextValue $ws.Cells.Item(16, 5) '  +2.73%  '

Set-TextValue $ws.Cells.Item(17, 4) '67.039.09'
Set-TextValue $ws.Cells.Item(17, 5) '  +1.24%  '

Set-TextValue $ws.Cells.Item(18, 4) '7.20'
Set-TextValue $ws.Cells.Item(18, 5) '  -0.87%  '

Set-TextValue $ws.Cells.Item(19, 4) '3.115.06'
Set-TextValue $ws.Cells.Item(19, 5) '  +2.83%  '

Set-TextValue $ws.Cells.Item(20, 4) '16.21'
Set-TextValue $ws.Cells.Item(20, 5) '  -2.01%  '

Set-TextValue $ws.Cells.Item(21, 4) '486.05'
Set-TextValue $ws.Cells.Item(21, 5) '  +2.39%  '

Set-TextValue $ws.Cells.Item(22, 4) '0.719'
Set-TextValue $ws.Cells.Item(22, 5) '  +1.42%  '

Set-TextValue $ws.Cells.Item(23, 4) '7.62'
Set-TextValue $ws.Cells.Item(23, 5) '  +1.84%  '

Set-TextValue $ws.Cells.Item(24, 4) '84.45'
Set-TextValue $ws.Cells.Item(24, 5) '  +1.40%  '

Set-TextValue $ws.Cells.Item(25, 4) '13.35'
Set-TextValue $ws.Cells.Item(25, 5) '  +3.79%  '

Set-TextValue $ws.Cells.Item(26, 5) '  +3.90%  '

Set-TextValue $ws.Cells.Item(27, 4) '10.08'
Set-TextValue $ws.Cells.Item(27, 5) '  +0.02%  '

Set-TextValue $ws.Cells.Item(28, 5) '  -0.03%  '

Set-TextValue $ws.Cells.Item(29, 4) '8.02'
Set-TextValue $ws.Cells.Item(29, 5) '  -2.41%  '

Set-TextValue $ws.Cells.Item(30, 4) '2.41'
Set-TextValue $ws.Cells.Item(30, 5) '  -1.68%  '

Set-TextValue $ws.Cells.Item(31, 5) '  +2.29%  '

Set-TextValue $ws.Cells.Item(32, 4) '28.88'
Set-TextValue $ws.Cells.Item(32, 5) '  +2.86%  '

Set-TextValue $ws.Cells.Item(33, 5) '  +0.49%  '

Set-TextValue $ws.Cells.Item(34, 4) '0.115'
Set-TextValue $ws.Cells.Item(34, 5) '  -2.38%  '

Set-TextValue $ws.Cells.Item(35, 5) '  +0.08%  '

Set-TextValue $ws.Cells.Item(36, 5) '  +1.07%  '

Set-TextValue $ws.Cells.Item(37, 4) '0.988'
Set-TextValue $ws.Cells.Item(37, 5) '  -0.32%  '

Set-TextValue $ws.Cells.Item(38, 4) '47.80'
Set-TextValue $ws.Cells.Item(38, 5) '  -0.63%  '

Set-TextValue $ws.Cells.Item(39, 5) '  +2.67%  '

Set-TextValue $ws.Cells.Item(40, 4) '50.18'
Set-TextValue $ws.Cells.Item(40, 5) '  +1.05%  '

Set-TextValue $ws.Cells.Item(41, 5) '  +1.54%  '

Set-TextValue $ws.Cells.Item(42, 5) '  +1.58%  '

Set-TextValue $ws.Cells.Item(43, 5) '  +0.42%  '

Set-TextValue $ws.Cells.Item(44, 4) '2.82'
Set-TextValue $ws.Cells.Item(44, 5) '  -0.98%  '

Set-TextValue $ws.Cells.Item(45, 4) '2.846.18'
Set-TextValue $ws.Cells.Item(45, 5) '  +4.45%  '

Set-TextValue $ws.Cells.Item(46, 5) '  -0.45%  '

Set-TextValue $ws.Cells.Item(47, 4) '383.39'
Set-TextValue $ws.Cells.Item(47, 5) '  -0.16%  '

Set-TextValue $ws.Cells.Item(48, 4) '137.54'
Set-TextValue $ws.Cells.Item(48, 5) '  +2.13%  '

Set-TextValue $ws.Cells.Item(50, 4) '25.25'
Set-TextValue $ws.Cells.Item(50, 5) '  +2.61%  '

Set-TextValue $ws.Cells.Item(51, 5) '  -0.25%  '
